# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 122
$ws1.Range("F3").Value = 2139
$ws1.Range("F4").Value = 24
$ws1.Range("F5").Value = 11167
$ws1.Range("F7").Value = 171
$ws1.Range("F8").Value = 308
$ws1.Range("F10").Value = 11073
$ws1.Range("F12").Value = 1141
$ws1.Range("F14").Value = 1724
$ws1.Range("F15").Value = 5552
$ws1.Range("F16").Value = 91
$ws1.Range("F17").Value = 3436

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 122
$ws4.Range("F3").Value = 2139
$ws4.Range("F5").Value = 24
$ws4.Range("F7").Value = 11167
$ws4.Range("F9").Value = 171
$ws4.Range("F10").Value = 308
$ws4.Range("F12").Value = 11073
$ws4.Range("F14").Value = 1141
$ws4.Range("F16").Value = 1724
$ws4.Range("F17").Value = 5552
$ws4.Range("F18").Value = 91
$ws4.Range("F19").Value = 3436
